$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "an email: nrzmhi@uni-wuerzburg.de" (sz 36 header block) - merge runs, drop proofErr
ReplaceText "an email: nrzmhi@uni-wuerzburg.de" "an email: nrzmhi@uni-wuerzburg.de"

# 2. "Wir möchten Sie höflich bitten..." sentence - merge runs, drop proofErr
ReplaceText "Wir möchten Sie höflich bitten, nach Erhalt des Befundes eine Empfangsbestätigung an die Nummer: 0931-31 87281 oder an die Email Adresse: nrzmhi@uni-wuerzburg.de zu senden." "Wir möchten Sie höflich bitten, nach Erhalt des Befundes eine Empfangsbestätigung an die Nummer: 0931-31 87281 oder an die Email Adresse: nrzmhi@uni-wuerzburg.de zu senden."
